$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite the previously-blank last row (row 23) with the new
# "Asphalt shingle roofing" entry (no row shift needed here).
$ws.Range("A23").ClearFormats()
$ws.Range("A23").Value = "Asphalt shingle roofing"
$ws.Range("B23").Value = 20

# Insert "Glazing: operable window" before the old row 22 (EPDM roofing).
$ws.Rows.Item(22).Insert()
$ws.Range("A22").ClearFormats()
$ws.Range("A22").Value = "Glazing: operable window"
$ws.Range("B22").Value = 30

# Insert "Brick: wood framing" before the old row 14 (Insulated metal panel).
$ws.Rows.Item(14).Insert()
$ws.Range("A14").ClearFormats()
$ws.Range("A14").Value = "Brick: wood framing"
$ws.Range("B14").Value = 35

# Insert four new rows before the old row 10 (Curtain wall: steel spandrel),
# in reverse order so the final reading order is:
#   Floor framing, Sub-flooring, Roof framing, Roof decking
$ws.Rows.Item(10).Insert()
$ws.Range("A10").ClearFormats()
$ws.Range("A10").Value = "Roof decking"
$ws.Range("B10").Value = 60

$ws.Rows.Item(10).Insert()
$ws.Range("A10").ClearFormats()
$ws.Range("A10").Value = "Roof framing"
$ws.Range("B10").Value = 60

$ws.Rows.Item(10).Insert()
$ws.Range("A10").ClearFormats()
$ws.Range("A10").Value = "Sub-flooring"
$ws.Range("B10").Value = 60

$ws.Rows.Item(10).Insert()
$ws.Range("A10").ClearFormats()
$ws.Range("A10").Value = "Floor framing"
$ws.Range("B10").Value = 60

# Insert "Concrete footing" before the old row 3 (Wall foundation).
$ws.Rows.Item(3).Insert()
$ws.Range("A3").ClearFormats()
$ws.Range("A3").Value = "Concrete footing"
$ws.Range("B3").Value = 60

# Update the selected cell to match the post-edit workbook state.
$ws.Range("E13").Select()
